$d = $word.ActiveDocument

# Locate the exact sentence that needs to be highlighted in yellow:
# "Beaudei also wants to have an image-slider using JQuery on this page
#  to show the latest events, promotions or news."
# This is the second "Beaudei ..." sentence in the "Home" bullet point
# paragraph (the first "Beaudei wants this page ..." sentence must stay
# untouched).
$r = $d.Content
$found = $r.Find.Execute("Beaudei also wants to have an image-slider using JQuery on this page to show the latest events, promotions or news.", `
    $false, $true, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Apply a yellow highlight to the found range only (wdYellow = 7).
    # Using Range.Font.HighlightColorIndex applies the highlight precisely
    # to the matched run span (splitting runs as needed), unlike setting
    # HighlightColorIndex directly on the Range object.
    $r.Font.HighlightColorIndex = 7
}
